$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('C5').Value = 53.05
$ws.Range('D5').Value = 53.05
$ws.Range('E5').Value = 53.05
$ws.Range('F5').Value = 53.05
$ws.Range('C6').Value = 1.03
$ws.Range('D6').Value = 0.79
$ws.Range('E6').Value = 1.35
$ws.Range('F6').Value = 1.57
$ws.Range('C7').Value = 19.26
$ws.Range('D7').Value = 20.9
$ws.Range('E7').Value = 17.73
$ws.Range('F7').Value = 13.39
$ws.Range('C8').Value = 26.66
$ws.Range('D8').Value = 25.26
$ws.Range('E8').Value = 27.87
$ws.Range('F8').Value = 31.99
$ws.Range('C16').Value = 4
$ws.Range('C17').Value = 1
$ws.Range('C18').Value = 6
$ws.Range('C28').Value = 0.36
$ws.Range('D28').Value = 0.64
$ws.Range('C29').Value = 0.55
$ws.Range('D29').Value = 0.45
$ws.Range('B38').Value = 2011
$ws.Range('C38').Value = 2.5725
$ws.Range('D38').Value = 234.16
$ws.Range('B39').Value = 2012
$ws.Range('C39').Value = 3.1618
$ws.Range('D39').Value = 277.4
$ws.Range('B40').Value = 2013
$ws.Range('C40').Value = 3.8915
$ws.Range('D40').Value = 309.99
$ws.Range('B41').Value = 2014
$ws.Range('C41').Value = 4.4698
$ws.Range('D41').Value = 318.74
$ws.Range('B42').Value = 2015
$ws.Range('C42').Value = 4.9126
$ws.Range('D42').Value = 461.98
$ws.Range('C50').Value = 'NA'
$ws.Range('D50').Value = 0.26
$ws.Range('B51').Value = 'Non-Promoter'
$ws.Range('C51').Value = 4.91
$ws.Range('D51').Value = 0.36
$ws.Range('C59').Value = 'Dividend (LHS)'
$ws.Range('D59').Value = 'EPS (LHS)'
$ws.Range('E59').Value = 'Payout (RHS)'
$ws.Range('B60').Value = 2013
$ws.Range('C60').Value = 2.1
$ws.Range('D60').Value = 8.17
$ws.Range('E60').Value = 0.3
$ws.Range('B61').Value = 2014
$ws.Range('C61').Value = 2.1
$ws.Range('D61').Value = 8.35
$ws.Range('E61').Value = 0.29
$ws.Range('B62').Value = 2015
$ws.Range('C62').Value = 2.3
$ws.Range('D62').Value = 9.27
$ws.Range('E62').Value = 0.29
$ws.Range('B71').Value = 'Dividend'
$ws.Range('C71').Value = 'EPS'
$ws.Range('D71').Value = 'Payout'
$ws.Range('B72').Value = 2.3
$ws.Range('C72').Value = 9.27
$ws.Range('D72').Value = 0.29
$ws.Range('B73').Value = 1
$ws.Range('C73').Value = 3.9
$ws.Range('D73').Value = 0.3
$ws.Range('B74').Value = 3.2
$ws.Range('C74').Value = 21.12
$ws.Range('D74').Value = 0.18
$ws.Range('C82').Value = 2015
$ws.Range('D82').Value = 2014
$ws.Range('B83').Value = 'Audit'
$ws.Range('C83').Value = 1.22
$ws.Range('D83').Value = 1.22
$ws.Range('B84').Value = 'Audit-Related'
$ws.Range('C84').Value = 0.41
$ws.Range('D84').Value = 0.38
$ws.Range('B85').Value = 'Non Audit'
$ws.Range('C85').Value = 0.29
$ws.Range('D85').Value = 0.33
$ws.Range('C96').Value = 'FY 12/13'
$ws.Range('D96').Value = 'FY 13/14'
$ws.Range('E96').Value = 'FY 14/15'
$ws.Range('B97').Value = 'Audit Fee'
$ws.Range('C97').Value = 1.02
$ws.Range('D97').Value = 1.22
$ws.Range('E97').Value = 1.22
$ws.Range('B98').Value = 'Audit Related Fee'
$ws.Range('C98').Value = 0.362
$ws.Range('D98').Value = 0.38
$ws.Range('E98').Value = 0.41
$ws.Range('B99').Value = 'Non Audit Fee '
$ws.Range('C99').Value = 0.28
$ws.Range('D99').Value = 0.33
$ws.Range('E99').Value = 0.29
$ws.Range('C109').Value = 'ED Remuneration'
$ws.Range('D109').Value = 'Indexed TSR'
$ws.Range('E109').Value = 'Net Profit'
$ws.Range('B110').Value = 'FY 14/15'
$ws.Range('C110').Value = 4.9126
$ws.Range('D110').Value = 461.98
$ws.Range('E110').Value = 823.07
$ws.Range('B111').Value = 'FY 13/14'
$ws.Range('C111').Value = 4.4698
$ws.Range('D111').Value = 318.74
$ws.Range('E111').Value = 741.14
$ws.Range('B112').Value = 'FY 12/13'
$ws.Range('C112').Value = 3.8915
$ws.Range('D112').Value = 309.99
$ws.Range('E112').Value = 725.18
$ws.Range('B113').Value = 'FY 11/12'
$ws.Range('C113').Value = 3.1618
$ws.Range('D113').Value = 277.4
$ws.Range('E113').Value = 600.16
$ws.Range('B114').Value = 'FY 10/11'
$ws.Range('C114').Value = 2.5725
$ws.Range('D114').Value = 234.16
$ws.Range('E114').Value = 432.61
$ws.Range('B115').Value = 'FY 9/10'
$ws.Range('C115').Value = 0
$ws.Range('D115').Value = 0
$ws.Range('E115').Value = 0
$ws.Range('C122').Value = ' '
$ws.Range('D122').Value = 'Commission'
$ws.Range('C123').Value = 'Promoter NED'
$ws.Range('D123').Value = 24.96
$ws.Range('C124').Value = 'Independent Directors'
$ws.Range('D124').Value = 34.49
$ws.Range('C125').Value = 'Other NEDs'
$ws.Range('D125').Value = 0
$ws.Range('C133').Value = ' '
$ws.Range('D133').Value = 'Total Commission'
$ws.Range('C134').Value = 'FY 10/11'
$ws.Range('D134').Value = 0
$ws.Range('C135').Value = 'FY 11/12'
$ws.Range('D135').Value = 0
$ws.Range('C136').Value = 'FY 12/13'
$ws.Range('D136').Value = 0
$ws.Range('C137').Value = 'FY 13/14'
$ws.Range('D137').Value = 0
$ws.Range('C138').Value = 'FY 14/15'
$ws.Range('D138').Value = 30.68
$ws.Range('C153').Value = 'Company'
$ws.Range('D153').Value = 'S&P CNX Nifty'
$ws.Range('E153').Value = 'CNX Finance'
$ws.Range('B154').Value = 'Today'
$ws.Range('C154').Value = 90
$ws.Range('D154').Value = 234
$ws.Range('E154').Value = 354
$ws.Range('B155').Value = '1Y'
$ws.Range('C155').Value = 100
$ws.Range('D155').Value = 213
$ws.Range('E155').Value = 298
$ws.Range('B156').Value = '3Y'
$ws.Range('C156').Value = 130
$ws.Range('D156').Value = 235
$ws.Range('E156').Value = 675
$ws.Range('B157').Value = '5Y'
$ws.Range('C157').Value = 140
$ws.Range('D157').Value = 987
$ws.Range('E157').Value = 283
$ws.Range('B483').Value = ' '
$ws.Range('C483').Value = 'Existing Borrowing'
$ws.Range('D483').Value = 'Unavailed borrowing limit'
$ws.Range('E483').Value = 'Proposed increase'
$ws.Range('B484').Value = 'Dec''15'
$ws.Range('C484').Value = 123
$ws.Range('D484').Value = 216
$ws.Range('B485').Value = 'Dec''14'
$ws.Range('C485').Value = 187
$ws.Range('D485').Value = 987
$ws.Range('E485').Value = 1009
$ws.Range('C493').Value = 'Dec''15'
$ws.Range('D493').Value = 23
$ws.Range('C494').Value = 'Dec''14'
$ws.Range('D494').Value = 24
$ws.Range('C495').Value = 'Dec''13'
$ws.Range('D495').Value = 26.5
$ws.Range('C496').Value = 'Dec''12'
$ws.Range('D496').Value = 21.87
$ws.Range('C504').Value = 'Dec''15'
$ws.Range('D504').Value = 12.34
$ws.Range('C505').Value = 'Dec''15'
$ws.Range('D505').Value = 21.21
$ws.Range('C506').Value = 'Dec''15'
$ws.Range('D506').Value = 12.34
$ws.Range('C507').Value = 'Dec''15'
$ws.Range('D507').Value = 64.32
$ws.Range('C515').Value = ' '
$ws.Range('D515').Value = 'CSR'
$ws.Range('E515').Value = 'CSR as % of Net Profit'
$ws.Range('C516').Value = 'FY 12/13'
$ws.Range('D516').Value = 213
$ws.Range('E516').Value = 28
$ws.Range('C517').Value = 'FY 13/14'
$ws.Range('D517').Value = 125
$ws.Range('E517').Value = 21
$ws.Range('C518').Value = 'FY 14/15'
$ws.Range('D518').Value = 123
$ws.Range('E518').Value = 23
$ws.Range('C531').Value = 'ED Remuneration'
$ws.Range('D531').Value = 'Indexed TSR'
$ws.Range('E531').Value = 'Net Profit'
$ws.Range('B532').Value = 'FY 14/15'
$ws.Range('C532').Value = 4.9126
$ws.Range('D532').Value = 461.98
$ws.Range('E532').Value = 823.07
$ws.Range('B533').Value = 'FY 13/14'
$ws.Range('C533').Value = 4.4698
$ws.Range('D533').Value = 318.74
$ws.Range('E533').Value = 741.14
$ws.Range('B534').Value = 'FY 12/13'
$ws.Range('C534').Value = 3.8915
$ws.Range('D534').Value = 309.99
$ws.Range('E534').Value = 725.18
$ws.Range('B535').Value = 'FY 11/12'
$ws.Range('C535').Value = 3.1618
$ws.Range('D535').Value = 277.4
$ws.Range('E535').Value = 600.16
$ws.Range('B536').Value = 'FY 10/11'
$ws.Range('C536').Value = 2.5725
$ws.Range('D536').Value = 234.16
$ws.Range('E536').Value = 432.61
$ws.Range('B537').Value = 'FY 9/10'
$ws.Range('C537').Value = 0
$ws.Range('D537').Value = 0
$ws.Range('E537').Value = 0
